$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tag matrix table currently lists 11 element/tag rows (row 1 is the header).
# A new "secure testversion" adds firewall1/firewall2 elements (and their
# to-firewall1 / to-firewall2 data-flow tags) into the existing matrix.
# Insert 4 blank rows (after row 4 "client") so the table grows from 11 to 15
# rows; new rows inherit the same style as the surrounding data rows.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Now (re)write column A for every data row (2-15) with the final tag order,
# and column B only where a value ("X") should remain.
$ws.Cells.Item(2,1).Value2  = "app"
$ws.Cells.Item(3,1).Value2  = "to-database"
$ws.Cells.Item(4,1).Value2  = "client"
$ws.Cells.Item(5,1).Value2  = "to-firewall1"
$ws.Cells.Item(6,1).Value2  = "database"
$ws.Cells.Item(7,1).Value2  = "firewall1"
$ws.Cells.Item(8,1).Value2  = "to-webapp"
$ws.Cells.Item(9,1).Value2  = "firewall2"
$ws.Cells.Item(10,1).Value2 = "to-app"
$ws.Cells.Item(11,1).Value2 = "webapp"
$ws.Cells.Item(12,1).Value2 = "to-firewall2"
$ws.Cells.Item(13,1).Value2 = "user-data"
$ws.Cells.Item(14,1).Value2 = "database"
$ws.Cells.Item(15,1).Value2 = "dmz"

# Clear any leftover column B values from the shift, then mark the
# "user-data" row (13) with an X, same as before the edit.
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r,2).Value2 = ""
}
$ws.Cells.Item(13,2).Value2 = "X"
